$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# TC_04 - Successful Admin Registration
$ws.Range("A16").Value = "TC_04"
$ws.Range("B16").Value = "Successful Admin Registration"
$ws.Range("C18").Value = "3.Type ""newadmin@nesto.com"" into ""//input[@placeholder='admin@nesto.com']"""
$ws.Range("C17").Value = "2.Type ""Admin User"" into ""//input[@placeholder='Enter Full Name']"""
$ws.Range("C16").Value = '1.Open URL "http://localhost:8080/signup"'
$ws.Range("C19").Value = "4.Type ""password123"" into ""//input[@placeholder='Create Password']"""
$ws.Range("C20").Value = "5.Type ""30"" into ""//input[@placeholder='Your Age']"""
$ws.Range("C21").Value = "6.Type ""9876543210"" into ""//input[@placeholder='Your Mobile Number']"""
$ws.Range("C22").Value = "7.Click on the ""Sign Up"" button ""//button[@type='submit']"""
$ws.Range("C23").Value = '8.Verify that the URL is "login"'

# TC_05 - Navigate back to Login
$ws.Range("A24").Value = "TC_05"
$ws.Range("B24").Value = "Navigate back to Login"
$ws.Range("C25").Value = "2.Click on ""Login here"" at ""//a[contains(text(), 'Login here')]"""
$ws.Range("C26").Value = '3.Verify that the URL is "login"'
$ws.Range("C24").Value = '1.Open URL "http://localhost:8080/signup"'

$ws.Range("C26").Select()
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1
